$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Stage") for the Level-4 module rows (38-52) previously held the
# text "M" (a stray/placeholder stage marker). Replace it with the actual
# numeric stage, matching the Level column pattern used throughout the sheet:
# rows 38-44 -> Stage 4, rows 45-52 -> Stage 5.
$ws.Range("E38:E44").Value = 4
$ws.Range("E45:E52").Value = 5

# Update the sheet's saved view/selection state.
$ws.Application.Goto($ws.Range("A26"), $true)
$ws.Range("G45:G52").Select()
